$wb = $excel.ActiveWorkbook

# The sheet "赵羽佳" (Zhao Yujia) is the 5th sheet (index 4, 0-based) / Worksheets.Item(5)
$ws = $wb.Worksheets.Item(5)

# --- New row 5 content -----------------------------------------------------
# A5: date text "2012.4.9" - must land in shared strings as plain text, not
# get auto-converted to a date serial. Force text format, write it, then
# clear the format again so the cell ends up with the default (no explicit)
# style, matching neighbouring date cells (A2/A3/A4).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2012.4.9"
$ws.Range("A5").ClearFormats()

# B5: the work-log note for this entry.
$workNote = "1.修改enemyball的运动方式函数，`n不同类型的球有不同的运行轨迹`n2.修改球的生命值与颜色变化。1-红、`n  2-橙、3-黄、4-绿、其他未白色。`n3.修改@property的习惯写法。@property成员变量名，变为_+名字"
$ws.Range("B5").Value = $workNote

# C5: the risk/issue note for this entry.
$riskNote = "球的血量现在是以颜色辨别。`n为了不让颜色太多而使人眼花，`n现在暂定有5种颜色。"
$ws.Range("C5").Value = $riskNote

# D5: numeric progress value.
$ws.Range("D5").Value = 3

# B5/C5 use the same wrap-text style as the rest of the table's note column.
$ws.Range("B5:C5").WrapText = $true

# Row 5 is a tall, wrapped row - match the authored row height.
$ws.Rows.Item(5).RowHeight = 81

# --- Column width tweak ------------------------------------------------
# Column C widens from 24 to 25 characters. Excel's ColumnWidth setter
# measures in character widths that get re-derived from pixels (+5/7 here),
# so back it out to land exactly on a stored width of 25.
$ws.Columns.Item(3).ColumnWidth = 25 - 5/7

# --- Active sheet / tab selection --------------------------------------
# The active tab moves from "刘云鹏" (sheet index 3) to "赵羽佳" (sheet index 4).
# Activating this sheet updates workbookView's activeTab and moves
# tabSelected to this sheet's view (off of 刘云鹏's).
$ws.Activate()

# Selection on the newly active sheet moves to B5 (single cell).
$ws.Range("B5").Select()
